# Apply the changes described by the diff:
#  - Add a new "travel" subcategory between "miscellaneous" and "stock"
#    (inserts a row so the investment/savings categories shift down by one,
#    and the new cell inherits the red "expense" formatting used by the
#    other expense categories above it)
#  - Fix the "other" income/expense cell (I8) to use the green "income"
#    font formatting instead of the red "expense" formatting
#  - Update the active cell selection to I8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing stock/bond/crypto/fund/real_estate/savings
# down by one row. The new row inherits formatting from the row above (row 19,
# which uses the red "expense" style), matching the target style for "travel".
$ws.Rows.Item(20).Insert()
$ws.Range("I20").Value = "travel"

# Fix I8 ("other") to use the green "income" font style (same style as I7)
# instead of the red "expense" style, by copying the format from I7.
$ws.Range("I7").Copy()
$ws.Range("I8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selected cell to I8
$ws.Range("I8").Select()
